$wb = $excel.ActiveWorkbook

$wsPurchase = $wb.Worksheets.Item("Purchase 22-23")
$wsSale     = $wb.Worksheets.Item("Sale 22-23")

# --- Data edit on "Purchase 22-23": row 31's bill was replaced ---
# Old: SLH/3591 / Shree Laxmi Lighting Hub / 16535, dated 44981
# New: 22230489 / Bhagyalaxmi Electricals / 31382, dated 44987
$wsPurchase.Range("B31").Value = 44987
$wsPurchase.Range("C31").Value = 22230489
$wsPurchase.Range("D31").Value = "Bhagyalaxmi Electricals"
$wsPurchase.Range("E31").Value = 31382
$wsPurchase.Range("F31").Formula = "=E31"
$wsPurchase.Range("F31").Font.Bold = $false

# Serial numbers for the following two bill blocks shift down by one
$wsPurchase.Range("A33").Value = 8
$wsPurchase.Range("A35").Value = 9

# --- View state updates ---
# Selection on the "Sale 22-23" sheet moves to G14, and it is no longer the
# active tab.
$wsSale.Range("G14").Select()

# "Purchase 22-23" becomes the active sheet, scrolled so row 19 is visible,
# with the selection on F32.
$wsPurchase.Activate()
$wsPurchase.Range("F32").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

# The workbook window is restored (no longer minimized).
$excel.WindowState = -4143
